$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update sheet name
$ws.Name = "Through 2022-08-21"

# Update header label in I1
$ws.Range("I1").Value = "2022 (through 08-21)"

# Update I9 value (August row)
$ws.Range("I9").Value = 123

# Update I14 (Total) value
$ws.Range("I14").Value = 1094
